$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 539120.5600000001
$ws.Range("I15").Value = 539120.5600000001
$ws.Range("K15").Value = 1617361.68
$ws.Range("M15").Value = -1617192.68
$ws.Range("H19").Value = 11407686
$ws.Range("I19").Value = 9392067
$ws.Range("K19").Value = 9392067
$ws.Range("M19").Value = -9391892
$ws.Range("H20").Value = 6000
$ws.Range("I20").Value = 6000
$ws.Range("K20").Value = 6000
$ws.Range("M20").Value = -5770
$ws.Range("H35").Value = 6000
$ws.Range("I35").Value = 6000
$ws.Range("K35").Value = 6000
$ws.Range("M35").Value = -5621
$ws.Range("H40").Value = 2799.2
$ws.Range("I40").Value = 2750
$ws.Range("J40").Value = 2996
$ws.Range("K40").Value = 2750
$ws.Range("L40").Value = 2996
$ws.Range("M40").Value = -2575
$ws.Range("N40").Value = -3346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1690.6522
$ws.Range("I45").Value = 1329.25
$ws.Range("J45").Value = 4100
$ws.Range("K45").Value = 1329.25
$ws.Range("L45").Value = 4100
$ws.Range("M45").Value = -952.25
$ws.Range("N45").Value = -4854
$ws.Range("H63").Value = 1482.0834
$ws.Range("I63").Value = 1383.5714
$ws.Range("J63").Value = 1620
$ws.Range("K63").Value = 1383.5714
$ws.Range("L63").Value = 1620
$ws.Range("M63").Value = -697.5714
$ws.Range("N63").Value = -2992
$ws.Range("H66").Value = 1482.0834
$ws.Range("I66").Value = 1383.5714
$ws.Range("J66").Value = 1620
$ws.Range("K66").Value = 6917.857
$ws.Range("L66").Value = 8100
$ws.Range("M66").Value = -3485.857
$ws.Range("N66").Value = -14964

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20346
$ws.Range("H20").Value = 5115.6665
$ws.Range("I20").Value = 6500
$ws.Range("J20").Value = 4838.8
$ws.Range("K20").Value = 6500
$ws.Range("L20").Value = 4838.8
$ws.Range("M20").Value = -6253
$ws.Range("N20").Value = -5332.8
$ws.Range("H82").Value = 8500
$ws.Range("I82").Value = 8500
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 8500
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -8117
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 8500
$ws.Range("I85").Value = 8500
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 8500
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -7174
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7001700
$ws.Range("I6").Value = 17500000
$ws.Range("K6").Value = 17500000
$ws.Range("M6").Value = -17499887
$ws.Range("H7").Value = 85.333336
$ws.Range("I7").Value = 77.8
$ws.Range("K7").Value = 77.8
$ws.Range("M7").Value = 35.2
$ws.Range("H17").Value = 16200
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 16200
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 16200
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -16548
$ws.Range("H25").Value = 3000
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 15000
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 20000
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = -21250
$ws.Range("H59").Value = 10634.667
$ws.Range("I59").Value = 5952
$ws.Range("K59").Value = 5952
$ws.Range("M59").Value = -4807
$ws.Range("H60").Value = 9879.115
$ws.Range("I60").Value = 6985.7144
$ws.Range("J60").Value = 10945.105
$ws.Range("K60").Value = 6985.7144
$ws.Range("L60").Value = 10945.105
$ws.Range("M60").Value = -6474.7144
$ws.Range("N60").Value = -11967.105
$ws.Range("H62").Value = 2766.2222
$ws.Range("I62").Value = 2860
$ws.Range("J62").Value = 2649
$ws.Range("K62").Value = 2860
$ws.Range("L62").Value = 2649
$ws.Range("M62").Value = -2236
$ws.Range("N62").Value = -3897
$ws.Range("H65").Value = 2766.2222
$ws.Range("I65").Value = 2860
$ws.Range("J65").Value = 2649
$ws.Range("K65").Value = 14300
$ws.Range("L65").Value = 13245
$ws.Range("M65").Value = -11180
$ws.Range("N65").Value = -19485
$ws.Range("H68").Value = 319999
$ws.Range("J68").Value = 319999
$ws.Range("L68").Value = 319999
$ws.Range("N68").Value = -321497
$ws.Range("H71").Value = 319999
$ws.Range("J71").Value = 319999
$ws.Range("L71").Value = 959997
$ws.Range("N71").Value = -967485
$ws.Range("H74").Value = 24072.5
$ws.Range("I74").Value = 21300
$ws.Range("J74").Value = 24996.666
$ws.Range("K74").Value = 21300
$ws.Range("L74").Value = 24996.666
$ws.Range("M74").Value = -20426
$ws.Range("N74").Value = -26744.666
$ws.Range("H77").Value = 24072.5
$ws.Range("I77").Value = 21300
$ws.Range("J77").Value = 24996.666
$ws.Range("K77").Value = 63900
$ws.Range("L77").Value = 74989.99800000001
$ws.Range("M77").Value = -59532
$ws.Range("N77").Value = -83725.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1000
$ws.Range("J31").Value = 1000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3576
$ws.Range("H35").Value = 1511.0769
$ws.Range("J35").Value = 1511.0769
$ws.Range("L35").Value = 4533.2307
$ws.Range("N35").Value = -5109.2307
$ws.Range("H49").Value = 3000
$ws.Range("J49").Value = 3000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9312
$ws.Range("H57").Value = 4093.2
$ws.Range("I57").Value = 766
$ws.Range("J57").Value = 4925
$ws.Range("K57").Value = 2298
$ws.Range("L57").Value = 14775
$ws.Range("M57").Value = -1739
$ws.Range("N57").Value = -15893
$ws.Range("H63").Value = 3017.3635
$ws.Range("I63").Value = 997.75
$ws.Range("J63").Value = 4171.4287
$ws.Range("K63").Value = 2993.25
$ws.Range("L63").Value = 12514.2861
$ws.Range("M63").Value = -2244.25
$ws.Range("N63").Value = -14012.2861
$ws.Range("H66").Value = 3017.3635
$ws.Range("I66").Value = 997.75
$ws.Range("J66").Value = 4171.4287
$ws.Range("K66").Value = 8979.75
$ws.Range("L66").Value = 37542.85830000001
$ws.Range("M66").Value = -5235.75
$ws.Range("N66").Value = -45030.85830000001
$ws.Range("H74").Value = 1875.3334
$ws.Range("I74").Value = 313
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 939
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = 122
$ws.Range("N74").Value = -17122
$ws.Range("H77").Value = 1875.3334
$ws.Range("I77").Value = 313
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 2817
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = 2487
$ws.Range("N77").Value = -55608
$ws.Range("H94").Value = 3746.1428
$ws.Range("I94").Value = 661.5
$ws.Range("K94").Value = 1984.5
$ws.Range("M94").Value = -1308.5
$ws.Range("H99").Value = 1505.2
$ws.Range("I99").Value = 881.5
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2644.5
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = -398.5
$ws.Range("N99").Value = -16492
$ws.Range("H102").Value = 5491.25
$ws.Range("I102").Value = 3800
$ws.Range("J102").Value = 6055
$ws.Range("K102").Value = 11400
$ws.Range("L102").Value = 18165
$ws.Range("M102").Value = -8966
$ws.Range("N102").Value = -23033
$ws.Range("H122").Value = 5641.4287
$ws.Range("J122").Value = 686.3684
$ws.Range("L122").Value = 6177.3156
$ws.Range("N122").Value = -11077.3156
$ws.Range("H138").Value = 2178.3125
$ws.Range("I138").Value = 1215
$ws.Range("J138").Value = 2315.9285
$ws.Range("K138").Value = 3645
$ws.Range("L138").Value = 6947.7855
$ws.Range("M138").Value = 1495
$ws.Range("N138").Value = -17227.7855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1750.8334
$ws.Range("I68").Value = 1746
$ws.Range("J68").Value = 1751.2727
$ws.Range("K68").Value = 1746
$ws.Range("L68").Value = 1751.2727
$ws.Range("M68").Value = -997
$ws.Range("N68").Value = -3249.2727
$ws.Range("H71").Value = 1750.8334
$ws.Range("I71").Value = 1746
$ws.Range("J71").Value = 1751.2727
$ws.Range("K71").Value = 8730
$ws.Range("L71").Value = 8756.363499999999
$ws.Range("M71").Value = -4986
$ws.Range("N71").Value = -16244.3635
$ws.Range("H93").Value = 5372.846
$ws.Range("I93").Value = 6978.4736
$ws.Range("J93").Value = 1014.7143
$ws.Range("K93").Value = 6978.4736
$ws.Range("L93").Value = 1014.7143
$ws.Range("M93").Value = -5730.4736
$ws.Range("N93").Value = -3510.7143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 15000
$ws.Range("J14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("N14").Value = -15336
